# Fix vid tomma rader i excellen på deltagare
$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("ekipage").Name = "Ekipage"
$wb.Worksheets.Item("klasser").Name = "Klasser"
$wb.Worksheets.Item("linförare").Name = "Linförare"
$wb.Worksheets.Item("voltigörer").Name = "Tävlande"

# Update selections (active cell) on relevant sheets
$wsEkipage = $wb.Worksheets.Item("Ekipage")
$wsEkipage.Activate()
$wsEkipage.Range("S22").Select()

$wsKlasser = $wb.Worksheets.Item("Klasser")
$wsKlasser.Activate()
$wsKlasser.Range("B9").Select()

$wsTavlande = $wb.Worksheets.Item("Tävlande")
$wsTavlande.Activate()
$wsTavlande.Range("D102").Select()

# Re-activate first sheet for tabSelected
$wsEkipage.Activate()
